# The "<id>p023r_1</id>" text in the <div> header block is split across
# three separate runs (the "<id>" / "</id>" tag markup in Courier New,
# 7f6000, plus the "p023r_1" value in a plain black run sandwiched
# between them). Re-typing that exact text over itself via Find/Replace
# collapses it into a single run, taking on the formatting of the first
# run it touches (the Courier New tag styling), exactly as in the
# target revision.
$d = $word.ActiveDocument

$needle = "<id>p023r_1</id>"

$rng = $d.Content.Find.Execute($needle, $false, $false, $false, $false, `
                                $false, $true, 1, $false, $needle, 2)
